# Applies the fr_ParentText Crisis Video Scripts.docx translation edits
# (French locale strings for the "amour bienveillant" -> "amour et
# gentillesse" meditation-script wording refresh).
$d = $word.ActiveDocument

# Disable smart-quote / smart-apostrophe autocorrect so the straight
# apostrophes in the replacement text are not silently converted to
# curly ones when the new text is inserted via the object model.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# 1) "La pause du jour s'appelle la méditation de l'amour bienveillant."
#    -> "La séance de méditation d'aujourd'hui porte le nom de méditation
#        sur l'amour et la gentillesse."
$r = $d.Content
$found = $r.Find.Execute("La pause du jour s'appelle la méditation de l'amour bienveillant.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Edit 1 found: $found"
if ($found) {
    $r.Text = "La séance de méditation d'aujourd'hui porte le nom de méditation sur l'amour et la gentillesse."
}

# 2) "Demandez-vous, “Quelle est mon expérience en ce moment précis?” "
#    -> "Posez--vous la question de savoir, “Quelle est mon expérience en
#        ce moment précis ?” "
$r = $d.Content
$found = $r.Find.Execute("Demandez-vous, “Quelle est mon expérience en ce moment précis?” ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Edit 2 found: $found"
if ($found) {
    $r.Text = "Posez--vous la question de savoir, “Quelle est mon expérience en ce moment précis ?” "
}

# 3) "...placer une main sur votre cœur ou poitrine. "
#    -> "...placer une main sur votre cœur ou votre poitrine. "
$r = $d.Content
$found = $r.Find.Execute("Connectez-vous à votre cœur de manière aimable et douce. Vous pouvez placer une main sur votre cœur ou poitrine. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Edit 3 found: $found"
if ($found) {
    $r.Text = "Connectez-vous à votre cœur de manière aimable et douce. Vous pouvez placer une main sur votre cœur ou votre poitrine. "
}

# 4) "...envoyer des pensées d'amour bienveillant à votre enfant..."
#    -> "...envoyer des pensées d'amour et de gentillesse à votre enfant..."
$r = $d.Content
$found = $r.Find.Execute("Si vous vous sentez à l'aise, vous pouvez aussi envoyer des pensées d'amour bienveillant à votre enfant, votre partenaire, votre famille, et toute autre personne qui vous est proche.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Edit 4 found: $found"
if ($found) {
    $r.Text = "Si vous vous sentez à l'aise, vous pouvez aussi envoyer des pensées d'amour et de gentillesse à votre enfant, votre partenaire, votre famille, et toute autre personne qui vous est proche."
}
